# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (summary) sheet,
#    populated with the single fund holding for that quarter.
# 2. Prepend a corresponding "2022-Q1" row to the "总计" sheet, shifting the
#    existing rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: add the "2022-Q1" sheet, positioned right before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the page margins used by the other quarterly sheets (0.75/0.75/1/1/0.5/0.5 in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row (B1:H1) - store as text like the source data
$newSheet.Range("B1:H1").NumberFormat = "@"
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row 2 - single fund holding for 2022-Q1
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "002938"
$newSheet.Range("C2").Value = "中银证券健康产业灵活配置混合"
$newSheet.Range("D2").Value = "1.82"
$newSheet.Range("E2").Value = "59.84"
$newSheet.Range("F2").Value = "2.67"
$newSheet.Range("G2").Value = "0.0486"
$newSheet.Range("H2").Value = 8
# Drop the explicit Text number-format again so B2:G2 fall back to the
# workbook default style, matching the other quarterly sheets.
$newSheet.Range("B2:G2").ClearFormats()

# Copy over the bordered/bold header + index-column look from "2021-Q2"
$ws2021Q2 = $wb.Worksheets.Item("2021-Q2")
$ws2021Q2.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$ws2021Q2.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: prepend the "2022-Q1" row to the "总计" sheet
# ---------------------------------------------------------------------------
# Re-resolve the sheet by name: worksheet handles are positional, and the
# "总计" sheet moved one slot to the right when the new sheet was inserted
# in front of it.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.05

# Give the new index cell (A2) the same bordered style as the other index cells
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Renumber the index column for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# Restore original active sheet/selection
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
